# Auto-generated edit script: adds 12 new 'Business Audit' domain rows (rows 53-64)
# to Sheet1, including one rich-text (partially bold) prompt cell in B53.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53: Financial Reporting
$ws.Cells.Item(53, 1).Value = 'Financial Reporting'
$ws.Cells.Item(53, 2).Value = 'The organization is preparing its financial statements for year-end. You discover inconsistent revenue recognition practices across different subsidiaries. Investigate how this might violate IFRS 15 or local GAAP rules. Propose a standardized policy and controls to ensure consistent application of revenue recognition principles.'
$ws.Cells.Item(53, 3).Value = 'Business Audit'
$ws.Cells.Item(53, 4).Value = 'IFRS (IFRS 15), GAAP'

# Row 54: Revenue & Receivables
$ws.Cells.Item(54, 1).Value = 'Revenue & Receivables'
$ws.Cells.Item(54, 2).Value = 'During the audit, you find that accounts receivable turnover has drastically slowed, yet revenue figures remain high. Assess the risk of misstated revenue, referencing ISA requirements for audit evidence. Recommend procedures (e.g., confirmation of receivables, cut-off testing) to validate revenue accuracy.'
$ws.Cells.Item(54, 3).Value = 'Business Audit'
$ws.Cells.Item(54, 4).Value = 'ISA (International Standards on Auditing), IFRS / GAAP (Revenue Recognition)'

# Row 55: Inventory Management
$ws.Cells.Item(55, 1).Value = 'Inventory Management'
$ws.Cells.Item(55, 2).Value = 'A manufacturing client uses different valuation methods (FIFO, Weighted Average) across various plants without a clear policy. Examine how this inconsistency might affect the financial statements under IFRS (IAS 2) or GAAP. Propose controls and a unified approach to ensure consistent inventory valuation.'
$ws.Cells.Item(55, 3).Value = 'Business Audit'
$ws.Cells.Item(55, 4).Value = 'IFRS (IAS 2), GAAP'

# Row 56: Internal Controls over Financial Reporting
$ws.Cells.Item(56, 1).Value = 'Internal Controls over Financial Reporting'
$ws.Cells.Item(56, 2).Value = 'You find that management hasn’t performed a formal risk assessment or control review in two years. Describe how this violates COSO Internal Control principles. Outline a plan to conduct a control self-assessment and document key controls to meet SOX 404 requirements.'
$ws.Cells.Item(56, 3).Value = 'Business Audit'
$ws.Cells.Item(56, 4).Value = 'COSO Internal Control – Integrated Framework, Sarbanes-Oxley (SOX)'

# Row 57: Enterprise Risk Management
$ws.Cells.Item(57, 1).Value = 'Enterprise Risk Management'
$ws.Cells.Item(57, 2).Value = 'The company recently expanded into new markets but has not updated its risk register. Reference COSO ERM to identify potential strategic and operational risks. Propose a process to integrate these risks into the existing ERM framework and regularly review them at the board level.'
$ws.Cells.Item(57, 3).Value = 'Business Audit'
$ws.Cells.Item(57, 4).Value = 'COSO Enterprise Risk Management (ERM)'

# Row 58: Segment Reporting
$ws.Cells.Item(58, 1).Value = 'Segment Reporting'
$ws.Cells.Item(58, 2).Value = 'Management consolidates financial information from multiple segments without a transparent allocation of shared costs. Analyze the potential misstatement under IFRS 8 or GAAP segment reporting rules. Recommend ways to improve segment disclosures and ensure compliance.'
$ws.Cells.Item(58, 3).Value = 'Business Audit'
$ws.Cells.Item(58, 4).Value = 'IFRS 8, GAAP (Segment Reporting)'

# Row 59: Control Environment & Culture
$ws.Cells.Item(59, 1).Value = 'Control Environment & Culture'
$ws.Cells.Item(59, 2).Value = 'You observe a weak ‘tone at the top,’ where management does not enforce policies consistently. Discuss the importance of a strong control environment per COSO and how leadership’s behavior impacts compliance with SOX internal control requirements.'
$ws.Cells.Item(59, 3).Value = 'Business Audit'
$ws.Cells.Item(59, 4).Value = 'COSO Internal Control, SOX (Control Environment)'

# Row 60: Compliance with Local Regulations
$ws.Cells.Item(60, 1).Value = 'Compliance with Local Regulations'
$ws.Cells.Item(60, 2).Value = 'A subsidiary in another country is subject to local reporting standards that differ from IFRS. Management has not reconciled these differences in the group financials. Outline the potential compliance risks and propose steps for consistent reporting across jurisdictions.'
$ws.Cells.Item(60, 3).Value = 'Business Audit'
$ws.Cells.Item(60, 4).Value = 'IFRS (Group Reporting), Local GAAP, COSO Internal Control'

# Row 61: Revenue Recognition
$ws.Cells.Item(61, 1).Value = 'Revenue Recognition'
$ws.Cells.Item(61, 2).Value = 'A multinational company applies inconsistent revenue recognition policies across its subsidiaries, resulting in different cut-off dates and partial accruals. Investigate how this might breach IFRS 15 (Revenue from Contracts with Customers). Outline the internal control gaps using COSO Internal Control principles, and discuss how management can strengthen SOX Section 404 compliance'
$ws.Cells.Item(61, 3).Value = 'Business Audit'
$ws.Cells.Item(61, 4).Value = 'IFRS 15 (Revenue Recognition), COSO Internal Control (Control Activities, Monitoring), SOX (Section 404 – Internal Controls over Financial Reporting)'

# Row 62: Financial Close & Reporting
$ws.Cells.Item(62, 1).Value = 'Financial Close & Reporting'
$ws.Cells.Item(62, 2).Value = 'The quarterly close process is rushed, causing errors in financial statements. Identify the risks to IFRS compliance and highlight how COSO Internal Control can guide the design of a more robust close process. Evaluate how the organization can reduce its SOX compliance risks by improving period-end controls and reconciliations.'
$ws.Cells.Item(62, 3).Value = 'Business Audit'
$ws.Cells.Item(62, 4).Value = 'IFRS,COSO Internal Control, SOX (Sections 302, 404)'

# Row 63: Impairment Testing
$ws.Cells.Item(63, 1).Value = 'Impairment Testing'
$ws.Cells.Item(63, 2).Value = 'Management has not conducted an annual impairment review on intangible assets. Discuss how this could violate IFRS (IAS 36) and propose an internal control approach based on COSO to ensure timely impairment testing. Also consider how this deficiency could affect the company’s SOX control environment.'
$ws.Cells.Item(63, 3).Value = 'Business Audit'
$ws.Cells.Item(63, 4).Value = 'IFRS (IAS 36), COSO Internal Control, SOX (Control Environment, 404)'

# Row 64: Forecasting & Budgeting
$ws.Cells.Item(64, 1).Value = 'Forecasting & Budgeting'
$ws.Cells.Item(64, 2).Value = 'The company’s budgeting process lacks scenario analysis, resulting in unrealistic forecasts. Show how COSO ERM can help integrate strategic risk assessments into budgeting. Discuss how unreliable forecasts could lead to misstatements under IFRS and potentially trigger SOX compliance issues if forecasts are used in financial reporting.'
$ws.Cells.Item(64, 3).Value = 'Business Audit'
$ws.Cells.Item(64, 4).Value = 'COSO ERM (Strategic Planning, Risk Appetite), IFRS,SOX (Sections 302, 404 if forecasts impact reported figures)'

# Apply partial bold rich-text formatting to B53 (IFRS 15 / GAAP emphasis).
# The leading run is left untouched (no explicit formatting call) so it keeps
# the plain/default run with no rPr, matching how Excel only emits rPr for
# runs that were explicitly touched by a formatting action.
$ws.Cells.Item(53, 2).Characters(192, 7).Font.Bold = $true
$ws.Cells.Item(53, 2).Characters(199, 10).Font.Bold = $false
$ws.Cells.Item(53, 2).Characters(209, 4).Font.Bold = $true
$ws.Cells.Item(53, 2).Characters(213, 118).Font.Bold = $false

# Move the active selection to reflect where the author's cursor ended up
$ws.Range("B67").Select()

